$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: D5 label text change; delete E5 (must happen before B2 so shared-string order matches)
$ws.Range("D5").Value = "Coste sin incluir la potencia (€/h)"
$ws.Range("E5").Clear()

# Row 2: change B2 label and C2 formula
$ws.Range("B2").Value = "COSTE TOTAL INCLUYENDO POTENCIA (€):"
$ws.Range("C2").Formula = "=SUM(D6:D29)+F2*F3/F1"

# Row 3: B3 style change (use same style as B2 - centered with border)
$ws.Range("B3").HorizontalAlignment = -4108

# Rows 6-29: change D formula, clear E column
for ($r = 6; $r -le 29; $r++) {
    $ws.Range("D$r").Formula = "=IF(I$r=`"X`", C$r*`$I`$2,C$r*`$I`$1)"
    $ws.Range("E$r").Clear()
}

# Column B width (ColumnWidth uses a slightly different unit than the stored
# OOXML "width"; 39.17 here round-trips to a stored width of 40)
$ws.Columns("B").ColumnWidth = 39.17

# Selection
$ws.Range("C3").Select()
